# Fix Training Data Issue (#48)
#
# Column BF ("Date") was populated with a malformed date string
# ("6-10-2012-13") for every data row. The NBA box-score data was off by
# one day because of how NBA stats were reported, so the correct value
# for this sheet (6-10-2012-13.xlsx) is "2013-06-10".
#
# NOTE: a plain `$cell.Value = "2013-06-10"` assignment gets silently
# re-interpreted by Excel as a *date* (since the string matches a date
# pattern), which would turn the cell into a numeric date serial instead
# of keeping it as literal text. To preserve the original "plain text"
# cell type (matching how the sheet was originally authored), we build
# the replacement value as a text formula result and paste it back in as
# a value, which keeps the cell a string rather than Excel's
# "looks-like-a-date" auto-conversion.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("BF2:BF31")
$rng.Formula = "=""2013-06-10"""
$rng.Copy($null) | Out-Null
$rng.PasteSpecial(-4163)  # xlPasteValues
